$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 285-296 ---
$ws = $wb.Worksheets.Item("PIR")

$pirRows = @(
    @("2026-02-06","10:03:44","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:03:48","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:03:53","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:03:58","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:04:03","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:04:07","10:00","Bathroom","Motion Detected","Active"),
    @("2026-02-06","10:04:15","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:04:20","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:04:25","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:04:30","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:04:35","10:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","10:04:40","10:00","Bathroom","No Motion","Inactive")
)

$startRow = 285
$endRow = $startRow + $pirRows.Count - 1
# Column A holds text dates ("2026-02-06"); force Text format so Excel
# doesn't auto-convert them to date serial numbers.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $r = $startRow + $i
    $row = $pirRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# --- Humidity sheet: append rows 182-191 ---
$ws = $wb.Worksheets.Item("Humidity")

$humidityRows = @(
    @("2026-02-06","10:03:45","10:00","Bathroom","69.5%","Active"),
    @("2026-02-06","10:03:55","10:00","Bathroom","69.5%","Active"),
    @("2026-02-06","10:04:05","10:00","Bathroom","69.3%","Active"),
    @("2026-02-06","10:04:10","10:00","Bathroom","69.2%","Active"),
    @("2026-02-06","10:04:15","10:00","Bathroom","69.3%","Active"),
    @("2026-02-06","10:04:20","10:00","Bathroom","69.3%","Active"),
    @("2026-02-06","10:04:25","10:00","Bathroom","69.3%","Active"),
    @("2026-02-06","10:04:31","10:00","Bathroom","67.9%","Active"),
    @("2026-02-06","10:04:35","10:00","Bathroom","69.4%","Active"),
    @("2026-02-06","10:04:41","10:00","Bathroom","69.4%","Active")
)

$startRow = 182
$endRow = $startRow + $humidityRows.Count - 1
# Column A holds text dates, column E holds percentage-looking text
# ("69.5%"); force Text format on both so Excel keeps them as literal
# strings instead of converting to a date serial / numeric percentage.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
$ws.Range("E$startRow`:E$endRow").NumberFormat = "@"
for ($i = 0; $i -lt $humidityRows.Count; $i++) {
    $r = $startRow + $i
    $row = $humidityRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# --- Temperature sheet: append rows 182-191 ---
$ws = $wb.Worksheets.Item("Temperature")

$temperatureRows = @(
    @("2026-02-06","10:03:46","10:00","Bathroom","27.8C","Active"),
    @("2026-02-06","10:03:56","10:00","Bathroom","27.8C","Active"),
    @("2026-02-06","10:04:06","10:00","Bathroom","27.8C","Active"),
    @("2026-02-06","10:04:11","10:00","Bathroom","27.7C","Active"),
    @("2026-02-06","10:04:16","10:00","Bathroom","27.8C","Active"),
    @("2026-02-06","10:04:21","10:00","Bathroom","27.7C","Active"),
    @("2026-02-06","10:04:26","10:00","Bathroom","27.6C","Active"),
    @("2026-02-06","10:04:31","10:00","Bathroom","27.7C","Active"),
    @("2026-02-06","10:04:36","10:00","Bathroom","27.7C","Active"),
    @("2026-02-06","10:04:41","10:00","Bathroom","27.7C","Active")
)

$startRow = 182
$endRow = $startRow + $temperatureRows.Count - 1
# Column A holds text dates ("2026-02-06"); force Text format so Excel
# doesn't auto-convert them to date serial numbers.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
for ($i = 0; $i -lt $temperatureRows.Count; $i++) {
    $r = $startRow + $i
    $row = $temperatureRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
